# Restored from revision #195ae0bf3a96a4be83d88e277669f05aa2008e54.
# Update the "Good Morning" return-value cell to "Good Morning111".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E8").Value = "Good Morning111"
